$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# These cells hold numbers stored as text throughout the sheet, so keep
# the written values text-formatted to match the existing column type.
$ws.Range("C2:E2").NumberFormat = "@"
$ws.Range("C5:E5").NumberFormat = "@"

# Swap the runs/balls/fours figures between row 2 and row 5.
$ws.Range("C2").Value = "7"
$ws.Range("D2").Value = "10"
$ws.Range("E2").Value = "0"

$ws.Range("C5").Value = "26"
$ws.Range("D5").Value = "27"
$ws.Range("E5").Value = "4"
